# ScorecardKnowledge: update template and import data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Cell values (A:code, B:name_km, C:name_en header + 4 data rows)
# ----------------------------------------------------------------------
$data = @(
  @("code",   "name_km", "name_en"),
  @("sk_001", "មេរៀនម៉ូឌុលទី ១- ការណែនាំអំពីគណនេយ្យភាពសង្គម (ISAF)", "Module 1- introduction to ISAF"),
  @("sk_002", "មេរៀនម៉ូឌុលទី ២៖ ការសម្របសម្រួលការផ្សព្វផ្សាយព័ត៌មានសម្រាប់ប្រជាពលរដ្ឋ (I4C)", "Module 2: Facilitating the information for citizen (I4C)"),
  @("sk_003", "មេរៀនម៉ូឌុលទី ៣៖ ការសម្របសម្រួលប័ណ្ណដាក់ពិន្ទុសហគមន៍ (CSC) និង ការវាយតម្លៃខ្លួនឯងដោយអ្នកផ្តល់សេវា (SSA)", "Module3: Facilitating community scorecard and service provider self-assessment"),
  @("sk_004", "មេរៀនម៉ូឌុលទី ៤៖ ការសម្របសម្រួលកិច្ចប្រជុំរួម និងការរៀបចំផែនការរួមស្តីពីគណនេយ្យភាពសង្គម (JAAP)", "Module 4: Facilitating the interface meeting and JAAP.")
)

for ($r = 1; $r -le 5; $r++) {
  $ws.Cells.Item($r, 1).Value2 = $data[$r-1][0]
  $ws.Cells.Item($r, 2).Value2 = $data[$r-1][1]
  $ws.Cells.Item($r, 3).Value2 = $data[$r-1][2]
}

# ----------------------------------------------------------------------
# 2) Formatting
# ----------------------------------------------------------------------

# -- Header row (row 1): bold Arial, centered horizontally, bottom aligned.
#    A1/B1 use the automatic (no explicit theme) font color; C1 (new
#    column) carries the theme text color.
$a1 = $ws.Cells.Item(1,1)
$a1.ClearFormats()
$a1.Font.Name = "Arial"
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4107     # xlBottom

$b1 = $ws.Cells.Item(1,2)
$b1.ClearFormats()
$b1.Font.Name = "Arial"
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108   # xlCenter
$b1.VerticalAlignment = -4107     # xlBottom

$c1 = $ws.Cells.Item(1,3)
$c1.Font.Name = "Arial"
$c1.Font.Bold = $true
$c1.Font.ThemeColor = 1
$c1.HorizontalAlignment = -4108   # xlCenter
$c1.VerticalAlignment = -4107     # xlBottom

# -- Column A (code) data rows.
# A2 (first imported row) keeps the automatic (no explicit) font color.
$a2 = $ws.Cells.Item(2,1)
$a2.ClearFormats()
$a2.Font.Name = "Arial"
$a2.VerticalAlignment = -4107     # xlBottom

# A3:A5 use the theme colored Arial font.
$a3a5 = $ws.Range("A3:A5")
$a3a5.Font.Name = "Arial"
$a3a5.Font.ThemeColor = 1
$a3a5.VerticalAlignment = -4107   # xlBottom

# -- Column B (Khmer name) rows 2-5: Khmer OS Battambang, size 11, theme
#    color, solid white fill, no wrap/shrink.
$kmRange = $ws.Range("B2:B5")
$kmRange.Font.Name = "Khmer OS Battambang"
$kmRange.Font.Size = 11
$kmRange.Font.ThemeColor = 1
$kmRange.Interior.Color = 16777215          # solid white fill
$kmRange.WrapText = $false
$kmRange.ShrinkToFit = $false

# -- Column C (English name) rows 2-5: Arial, size 11, theme color, white fill.
$enRange = $ws.Range("C2:C5")
$enRange.Font.Name = "Arial"
$enRange.Font.Size = 11
$enRange.Font.ThemeColor = 1
$enRange.Interior.Color = 16777215          # solid white fill

Write-Host "Done"
